$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The ParticipantsTab query (cell B2) was rewritten to traverse the
# participant->study relationship in the other direction, add
# genomic_info/diagnosis optional matches, and sort the collected sample
# ids before joining them.
$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE s.study_name in ["Childhood Cancer Data Initiative (CCDI): Free the Data: Open Sharing of Comprehensive Genomic Childhood Cancer Datasets (Kansas)"]
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$ws.Range("B2").Value = $newQuery

# The cell wraps text, so the taller query pushes row 2's height out.
$ws.Rows.Item(2).RowHeight = 299.25

# The workbook was last saved with C3 selected rather than A2.
$ws.Range("C3").Select()
